$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.462.95"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "'1.875.31"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").Value = "313.63"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").Value = "0.4795"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").Value = "0.3772"
$ws.Range("E8").Value = "  +3.12%  "
$ws.Range("D9").Value = "0.07401"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").Value = "0.9429"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "20.82"
$ws.Range("E11").Value = "  +6.80%  "
$ws.Range("D12").Value = "0.07883"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").Value = "'1.888.55"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").Value = "5.456"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "6.633"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "91.05"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").Value = "1.017"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "0.000008975"
$ws.Range("E18").Value = "  +4.16%  "
$ws.Range("D19").Value = "1.014"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "14.96"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "'27.484.79"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "5.157"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("D23").Value = "10.71"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "1.956"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").Value = "154.15"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").Value = "'18.60"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").Value = "2.023"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "116.21"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "5.013"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").Value = "0.08937"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "3.329"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "1.221"
$ws.Range("E32").Value = "  +5.14%  "
$ws.Range("D33").Value = "4.597"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "0.7509"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "2.701"
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").Value = "0.02073"
$ws.Range("E36").Value = "  +6.54%  "
$ws.Range("D37").Value = "1.123"
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("D38").Value = "0.05318"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "'3.000"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "0.5355"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "7.094"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").Value = "0.1532"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "8.435"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("D44").Value = "10.68"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "0.4852"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "1.015"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").Value = "1.666"
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("D48").Value = "103.46"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("D49").Value = "67.28"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("D50").Value = "0.06114"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "0.9017"
$ws.Range("E51").Value = "  +2.19%  "
